$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'80.202.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.44%  "
$ws.Range("D3").Value = "'3.224.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.73%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'212.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.87%  "
$ws.Range("D6").Value = "'641.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("E7").Value = "  +31.44%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.88%  "
$ws.Range("D10").Value = "'3.224.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.82%  "
$ws.Range("D11").Value = "'0.623"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +42.32%  "
$ws.Range("E12").Value = "  +43.58%  "
$ws.Range("E13").Value = "  +3.69%  "
$ws.Range("D14").Value = "'5.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.35%  "
$ws.Range("D15").Value = "'3.818.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.86%  "
$ws.Range("D16").Value = "'32.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +14.45%  "
$ws.Range("D17").Value = "'80.032.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.27%  "
$ws.Range("D18").Value = "'3.220.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.95%  "
$ws.Range("D19").Value = "'14.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.32%  "
$ws.Range("D20").Value = "'3.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +28.75%  "
$ws.Range("E21").Value = "  +6.18%  "
$ws.Range("D22").Value = "'449.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +18.96%  "
$ws.Range("D23").Value = "'5.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +23.59%  "
$ws.Range("D24").Value = "'4.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.69%  "
$ws.Range("D25").Value = "'78.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.99%  "
$ws.Range("D26").Value = "'3.366.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.06%  "
$ws.Range("D27").Value = "'11.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.92%  "
$ws.Range("D28").Value = "'0.0000128"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +18.89%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +13.34%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'568.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.90%  "
$ws.Range("D33").Value = "'1.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.98%  "
$ws.Range("D34").Value = "'0.158"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +28.51%  "
$ws.Range("E35").Value = "  +7.94%  "
$ws.Range("D36").Value = "'23.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.42%  "
$ws.Range("E37").Value = "  +21.01%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +11.16%  "
$ws.Range("E40").Value = "  +15.19%  "
$ws.Range("D41").Value = "'164.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").Value = "'193.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  +13.44%  "
$ws.Range("D46").Value = "'2.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +15.12%  "
$ws.Range("E47").Value = "  +9.61%  "
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("D49").Value = "'4.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.07%  "
$ws.Range("D50").Value = "'43.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.97%  "
$ws.Range("D51").Value = "'0.652"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.53%  "
